$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add D1=3 and E1=4, matching style of C1 ---
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Bulk update B2:C67 (existing columns, new values) and D2:E67 (new columns) ---
$bc = New-Object "object[,]" 66,2
$de = New-Object "object[,]" 66,2

$bc[0,0] = -0.3168750346129532; $bc[0,1] = -0.3036867844311208
$de[0,0] = -0.2912284507030619; $de[0,1] = -0.2798558518241587
$bc[1,0] = 0.2264838909210563; $bc[1,1] = 0.226016869225491
$de[1,0] = 0.2265239587630224; $de[1,1] = 0.2273824790628174
$bc[2,0] = 0.1460605506844183; $bc[2,1] = 0.1547881991105919
$de[2,0] = 0.1638389968063508; $de[2,1] = 0.1725441346294188
$bc[3,0] = -0.0746019147878336; $bc[3,1] = -0.06892048399615523
$de[3,0] = -0.06362651820170304; $de[3,1] = -0.05903094877233309
$bc[4,0] = 0.1505426926751276; $bc[4,1] = 0.1588074336107037
$de[4,0] = 0.1662061879154078; $de[4,1] = 0.1724384878066593
$bc[5,0] = -0.4549700821850928; $bc[5,1] = -0.4459064887304521
$de[5,0] = -0.4366836754820621; $de[5,1] = -0.4277497093504777
$bc[6,0] = -0.2681618823914937; $bc[6,1] = -0.2533299822585643
$de[6,0] = -0.2388396335645593; $de[6,1] = -0.2251927607453582
$bc[7,0] = -0.3886549536128331; $bc[7,1] = -0.3869446842432583
$de[7,0] = -0.3833715160064274; $de[7,1] = -0.3786728188839405
$bc[8,0] = 0.3597580557313836; $bc[8,1] = 0.3689886109514408
$de[8,0] = 0.3772256881679152; $de[8,1] = 0.3840900581502769
$bc[9,0] = -0.247218855270827; $bc[9,1] = -0.2375676495630193
$de[9,0] = -0.2282243807386613; $de[9,1] = -0.2195786255325146
$bc[10,0] = -0.07745622716668545; $bc[10,1] = -0.06821648309498324
$de[10,0] = -0.06074958616600853; $de[10,1] = -0.05507102361751429
$bc[11,0] = -0.007112783682018264; $bc[11,1] = -0.009372778487127869
$de[11,0] = -0.01092949463124122; $de[11,1] = -0.01220303303378092
$bc[12,0] = 0.08881280458427744; $bc[12,1] = 0.09835750693104157
$de[12,0] = 0.1050442851898987; $de[12,1] = 0.1091322560133212
$bc[13,0] = 0.03002055446965043; $bc[13,1] = 0.0367282424715893
$de[13,0] = 0.03978028801853266; $de[13,1] = 0.03975663775094221
$bc[14,0] = 0.4073638009776918; $bc[14,1] = 0.4193125732117537
$de[14,0] = 0.426494213323987; $de[14,1] = 0.4295429358332071
$bc[15,0] = 0.6533922235476113; $bc[15,1] = 0.6483189581467216
$de[15,0] = 0.6403126297509089; $de[15,1] = 0.629865925182356
$bc[16,0] = 0.08345790035784724; $bc[16,1] = 0.05431395812032949
$de[16,0] = 0.02863334878819883; $de[16,1] = 0.005931364926652402
$bc[17,0] = 0.4219553194489087; $bc[17,1] = 0.4194873667791191
$de[17,0] = 0.4162318407594452; $de[17,1] = 0.4120570630130523
$bc[18,0] = 0.3503954740367983; $bc[18,1] = 0.3597478966727853
$de[18,0] = 0.3642107678940616; $de[18,1] = 0.3645378298914858
$bc[19,0] = 0.5333769738622589; $bc[19,1] = 0.5589686541212187
$de[19,0] = 0.5774431289596272; $de[19,1] = 0.5896719888669699
$bc[20,0] = 0.4218904866053518; $bc[20,1] = 0.4224199210089787
$de[20,0] = 0.4197146048390447; $de[20,1] = 0.4142825155653218
$bc[21,0] = -0.03488353509137355; $bc[21,1] = -0.02969368961751345
$de[21,0] = -0.02798139298564301; $de[21,1] = -0.02914775554181376
$bc[22,0] = 4.302167219374502; $bc[22,1] = 4.2964596053585
$de[22,0] = 4.232107234508156; $de[22,1] = 4.124280750287866
$bc[23,0] = 0.4983464021535683; $bc[23,1] = 0.4655444954022276
$de[23,0] = 0.4376316625048661; $de[23,1] = 0.41386762509255
$bc[24,0] = 0.3886796852984516; $bc[24,1] = 0.3692330055649125
$de[24,0] = 0.3501783215005168; $de[24,1] = 0.3318439063651439
$bc[25,0] = 0.2815464012767376; $bc[25,1] = 0.255308055126365
$de[25,0] = 0.2308337605734169; $de[25,1] = 0.2083636085980171
$bc[26,0] = 1.016151521223956; $bc[26,1] = 0.985866018235656
$de[26,0] = 0.9593702789342871; $de[26,1] = 0.9358612898752613
$bc[27,0] = 5.710310430599527; $bc[27,1] = 5.244289870569874
$de[27,0] = 4.814977541712624; $de[27,1] = 4.422688532356388
$bc[28,0] = 0.97845013991829; $bc[28,1] = 0.9343558192232184
$de[28,0] = 0.8954986204524927; $de[28,1] = 0.861280387428689
$bc[29,0] = -0.2391845546693665; $bc[29,1] = -0.2778607738411689
$de[29,0] = -0.3099547290097538; $de[29,1] = -0.3360917307501277
$bc[30,0] = 0.7799172826702934; $bc[30,1] = 0.7490073304897067
$de[30,0] = 0.7218855701460607; $de[30,1] = 0.6979985672338453
$bc[31,0] = 0.9097394503583777; $bc[31,1] = 0.885522784570417
$de[31,0] = 0.8649455041205639; $de[31,1] = 0.8472515887110703
$bc[32,0] = -0.6825464241570107; $bc[32,1] = -0.7030225518903331
$de[32,0] = -0.7204436003823272; $de[32,1] = -0.7349629998652148
$bc[33,0] = 0.7997055999329813; $bc[33,1] = 0.7996699859181595
$de[33,0] = 0.8004517126802467; $de[33,1] = 0.8012738614048697
$bc[34,0] = 0.7513764998381978; $bc[34,1] = 0.7480588837663376
$de[34,0] = 0.7468273871562174; $de[34,1] = 0.7466453032739981
$bc[35,0] = 0.7241841030854419; $bc[35,1] = 0.7194746032728944
$de[35,0] = 0.717039929643372; $de[35,1] = 0.7158423614551789
$bc[36,0] = 0.7194171909490212; $bc[36,1] = 0.7058849505409454
$de[36,0] = 0.6950575043410426; $de[36,1] = 0.6860194624807041
$bc[37,0] = 0.5651682208465996; $bc[37,1] = 0.5694553353364057
$de[37,0] = 0.574243812116516; $de[37,1] = 0.5787991313839743
$bc[38,0] = 0.7324010633978756; $bc[38,1] = 0.7387958335202249
$de[38,0] = 0.7449392446120547; $de[38,1] = 0.7502075611169975
$bc[39,0] = 0.5550718351305318; $bc[39,1] = 0.5510444178321311
$de[39,0] = 0.5489502942864384; $de[39,1] = 0.5478890744001382
$bc[40,0] = 0.6766021403910398; $bc[40,1] = 0.6647413438913339
$de[40,0] = 0.6559686245186239; $de[40,1] = 0.6492200940920476
$bc[41,0] = 0.7075626749135246; $bc[41,1] = 0.6990383841299922
$de[41,0] = 0.6927865310680269; $de[41,1] = 0.6878808432644599
$bc[42,0] = 0.6629724017749757; $bc[42,1] = 0.6635160940718288
$de[42,0] = 0.6654233746247672; $de[42,1] = 0.6677886944619366
$bc[43,0] = 0.6286851194373639; $bc[43,1] = 0.6310103914282831
$de[43,0] = 0.635397136946628; $de[43,1] = 0.6406915285061743
$bc[44,0] = -1.285318208976781; $bc[44,1] = -1.279920598506624
$de[44,0] = -1.274085652639843; $de[44,1] = -1.26798036835683
$bc[45,0] = -0.9979295053526136; $bc[45,1] = -0.9941628244536368
$de[45,0] = -0.9897895178064288; $de[45,1] = -0.9850697754450696
$bc[46,0] = -0.8923188071364571; $bc[46,1] = -0.8876851908265895
$de[46,0] = -0.8818629467764777; $de[46,1] = -0.8753429333227216
$bc[47,0] = -0.6580562947258826; $bc[47,1] = -0.6530785597146361
$de[47,0] = -0.6472912442600749; $de[47,1] = -0.641155311573986
$bc[48,0] = -0.05833912882412284; $bc[48,1] = -0.05576828822634386
$de[48,0] = -0.05283233324812737; $de[48,1] = -0.04998562428393243
$bc[49,0] = -0.8874537961833786; $bc[49,1] = -0.8808236328114679
$de[49,0] = -0.873410093434568; $de[49,1] = -0.8656399141082676
$bc[50,0] = -0.8874537961833786; $bc[50,1] = -0.8808236328114679
$de[50,0] = -0.873410093434568; $de[50,1] = -0.8656399141082676
$bc[51,0] = -1.112983025933588; $bc[51,1] = -1.114756251233469
$de[51,0] = -1.114283196105775; $de[51,1] = -1.112114272275124
$bc[52,0] = -0.2005019615151444; $bc[52,1] = -0.1941042360003489
$de[52,0] = -0.1874705415526502; $de[52,1] = -0.1810836286828947
$bc[53,0] = -1.018573519203073; $bc[53,1] = -1.013853905987674
$de[53,0] = -1.008695236059333; $de[53,1] = -1.003330689888567
$bc[54,0] = -0.9084862882514868; $bc[54,1] = -0.8976590279862675
$de[54,0] = -0.8880864840152791; $de[54,1] = -0.8797083979683353
$bc[55,0] = -0.936617848885817; $bc[55,1] = -0.9287350184708151
$de[55,0] = -0.9219440711288701; $de[55,1] = -0.9161488534922133
$bc[56,0] = -1.142249518574323; $bc[56,1] = -1.123567071138395
$de[56,0] = -1.107698600414897; $de[56,1] = -1.094281286249151
$bc[57,0] = -0.8661199415266272; $bc[57,1] = -0.8504870520082407
$de[57,0] = -0.835994533810961; $de[57,1] = -0.8227740022287811
$bc[58,0] = -0.5076263755722759; $bc[58,1] = -0.4903224098973188
$de[58,0] = -0.4751028818804094; $de[58,1] = -0.4619804466173932
$bc[59,0] = 0.3693360363631596; $bc[59,1] = 0.3718422085082914
$de[59,0] = 0.3747257306098983; $de[59,1] = 0.3774004402429267
$bc[60,0] = -1.216527512480593; $bc[60,1] = -1.202259018807272
$de[60,0] = -1.190177142068811; $de[60,1] = -1.179961085403063
$bc[61,0] = -0.7638644273104289; $bc[61,1] = -0.7343364021664399
$de[61,0] = -0.7073211609515856; $de[61,1] = -0.6829744799110793
$bc[62,0] = -0.9030195039197697; $bc[62,1] = -0.8974519606912046
$de[62,0] = -0.8908727515876083; $de[62,1] = -0.883748017420901
$bc[63,0] = -0.1241280864370546; $bc[63,1] = -0.106159846411747
$de[63,0] = -0.0903522647250089; $de[63,1] = -0.07684234004514823
$bc[64,0] = -0.8030625145166063; $bc[64,1] = -0.7851092787761388
$de[64,0] = -0.7705490646969968; $de[64,1] = -0.7589375968952176
$bc[65,0] = -0.7857817249430636; $bc[65,1] = -0.7594586154370843
$de[65,0] = -0.7387268444560317; $de[65,1] = -0.722746868067334

$ws.Range("B2:C67").Value2 = $bc
$ws.Range("D2:E67").Value2 = $de

Write-Output "done"